$wb = $excel.ActiveWorkbook

# -----------------------------------------------------------------------
# Sheet "Menu Mock": append new "CSM Advanced" rows 94-105 (debug
# diagnostics + quick-test controls) right after the existing row 93
# ("Debug Logging").
# -----------------------------------------------------------------------
$ws = $wb.Worksheets.Item("Menu Mock")

$ws.Cells.Item(94, 1).Value = "CSM Advanced"
$ws.Cells.Item(94, 2).Value = "Show Effective Values"
$ws.Cells.Item(94, 3).Value = "Toggle"
$ws.Cells.Item(94, 4).Value = "Off"
$ws.Cells.Item(94, 7).Value = "Show effective per-trigger values after presets and overrides"

$ws.Cells.Item(95, 1).Value = "CSM Advanced"
$ws.Cells.Item(95, 2).Value = "Last Trigger"
$ws.Cells.Item(95, 3).Value = "Arrow"
$ws.Cells.Item(95, 4).Value = "None"
$ws.Cells.Item(95, 6).Value = "LastTriggerProvider"
$ws.Cells.Item(95, 7).Value = "Last trigger attempt"

$ws.Cells.Item(96, 1).Value = "CSM Advanced"
$ws.Cells.Item(96, 2).Value = "Last Trigger Reason"
$ws.Cells.Item(96, 3).Value = "Arrow"
$ws.Cells.Item(96, 4).Value = "None"
$ws.Cells.Item(96, 6).Value = "LastTriggerReasonProvider"
$ws.Cells.Item(96, 7).Value = "Why the last trigger did or didn't fire"

$ws.Cells.Item(97, 1).Value = "CSM Advanced"
$ws.Cells.Item(97, 2).Value = "Quick Test Trigger"
$ws.Cells.Item(97, 3).Value = "Arrow"
$ws.Cells.Item(97, 4).Value = "Basic Kill"
$ws.Cells.Item(97, 5).Value = "Basic Kill | Critical Kill | Dismemberment | Decapitation | Parry | Last Enemy | Last Stand"
$ws.Cells.Item(97, 6).Value = "QuickTestTriggerProvider"
$ws.Cells.Item(97, 7).Value = "Which trigger to simulate"

$ws.Cells.Item(98, 1).Value = "CSM Advanced"
$ws.Cells.Item(98, 2).Value = "Quick Test Now"
$ws.Cells.Item(98, 3).Value = "Toggle"
$ws.Cells.Item(98, 4).Value = "Off"
$ws.Cells.Item(98, 7).Value = "Toggle to fire the selected trigger once"

$ws.Cells.Item(99, 1).Value = "CSM Advanced"
$ws.Cells.Item(99, 2).Value = "Effective: Basic Kill"
$ws.Cells.Item(99, 3).Value = "Arrow"
$ws.Cells.Item(99, 4).Value = "Off"
$ws.Cells.Item(99, 6).Value = "EffectiveBasicProvider"
$ws.Cells.Item(99, 7).Value = "Effective values for Basic Kill"

$ws.Cells.Item(100, 1).Value = "CSM Advanced"
$ws.Cells.Item(100, 2).Value = "Effective: Critical Kill"
$ws.Cells.Item(100, 3).Value = "Arrow"
$ws.Cells.Item(100, 4).Value = "Off"
$ws.Cells.Item(100, 6).Value = "EffectiveCriticalProvider"
$ws.Cells.Item(100, 7).Value = "Effective values for Critical Kill"

$ws.Cells.Item(101, 1).Value = "CSM Advanced"
$ws.Cells.Item(101, 2).Value = "Effective: Dismemberment"
$ws.Cells.Item(101, 3).Value = "Arrow"
$ws.Cells.Item(101, 4).Value = "Off"
$ws.Cells.Item(101, 6).Value = "EffectiveDismembermentProvider"
$ws.Cells.Item(101, 7).Value = "Effective values for Dismemberment"

$ws.Cells.Item(102, 1).Value = "CSM Advanced"
$ws.Cells.Item(102, 2).Value = "Effective: Decapitation"
$ws.Cells.Item(102, 3).Value = "Arrow"
$ws.Cells.Item(102, 4).Value = "Off"
$ws.Cells.Item(102, 6).Value = "EffectiveDecapitationProvider"
$ws.Cells.Item(102, 7).Value = "Effective values for Decapitation"

$ws.Cells.Item(103, 1).Value = "CSM Advanced"
$ws.Cells.Item(103, 2).Value = "Effective: Parry"
$ws.Cells.Item(103, 3).Value = "Arrow"
$ws.Cells.Item(103, 4).Value = "Off"
$ws.Cells.Item(103, 6).Value = "EffectiveParryProvider"
$ws.Cells.Item(103, 7).Value = "Effective values for Parry"

$ws.Cells.Item(104, 1).Value = "CSM Advanced"
$ws.Cells.Item(104, 2).Value = "Effective: Last Enemy"
$ws.Cells.Item(104, 3).Value = "Arrow"
$ws.Cells.Item(104, 4).Value = "Off"
$ws.Cells.Item(104, 6).Value = "EffectiveLastEnemyProvider"
$ws.Cells.Item(104, 7).Value = "Effective values for Last Enemy"

$ws.Cells.Item(105, 1).Value = "CSM Advanced"
$ws.Cells.Item(105, 2).Value = "Effective: Last Stand"
$ws.Cells.Item(105, 3).Value = "Arrow"
$ws.Cells.Item(105, 4).Value = "Off"
$ws.Cells.Item(105, 6).Value = "EffectiveLastStandProvider"
$ws.Cells.Item(105, 7).Value = "Effective values for Last Stand"

# -----------------------------------------------------------------------
# Sheet "Providers": add rows for the new value providers referenced
# above, keeping the existing alphabetical ordering of the list.
# -----------------------------------------------------------------------
$ws2 = $wb.Worksheets.Item("Providers")

# "Effective*Provider" entries sort alphabetically right before
# "GlobalSmoothingProvider" (currently row 18) -- insert 7 blank rows.
for ($i = 0; $i -lt 7; $i++) {
    $ws2.Rows.Item(18).Insert()
}

$ws2.Cells.Item(18, 1).Value = "EffectiveBasicProvider"
$ws2.Cells.Item(19, 1).Value = "EffectiveCriticalProvider"
$ws2.Cells.Item(20, 1).Value = "EffectiveDecapitationProvider"
$ws2.Cells.Item(21, 1).Value = "EffectiveDismembermentProvider"
$ws2.Cells.Item(22, 1).Value = "EffectiveLastEnemyProvider"
$ws2.Cells.Item(23, 1).Value = "EffectiveLastStandProvider"
$ws2.Cells.Item(24, 1).Value = "EffectiveParryProvider"

# "LastTrigger*Provider" entries sort alphabetically right before
# "MinEnemyGroupProvider" (now shifted down to row 30) -- insert 2 blank rows.
for ($i = 0; $i -lt 2; $i++) {
    $ws2.Rows.Item(30).Insert()
}

$ws2.Cells.Item(30, 1).Value = "LastTriggerProvider"
$ws2.Cells.Item(31, 1).Value = "LastTriggerReasonProvider"

# "QuickTestTriggerProvider" sorts alphabetically right before
# "SmoothingSpeedProvider" (now shifted down to row 34) -- insert 1 blank row.
$ws2.Rows.Item(34).Insert()

$ws2.Cells.Item(34, 1).Value = "QuickTestTriggerProvider"
$ws2.Cells.Item(34, 2).Value = "Basic Kill | Critical Kill | Dismemberment | Decapitation | Parry | Last Enemy | Last Stand"
